# Sample Project / Main.xlsx — "Rules" sheet
# Change B11 (Rule column, row "R40") so that its displayed value becomes
# the text "1" instead of "R40", keeping the cell's existing style/format.
#
# Assigning a plain numeric-looking string straight to .Value would make
# Excel auto-convert it to a real number (losing the shared-string/text
# nature of the cell, and changing its "t" type in the saved XML). To keep
# it a genuine text value we stage it in a scratch cell formatted as Text,
# then copy only the *value* (PasteSpecial values-only) onto B11 so the
# destination keeps its original number format / style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "1"

$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues

$scratch.Clear()
